$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 76
$ws.Range("H76").Value = 6388.4375
$ws.Range("I76").Value = 5378.5713
$ws.Range("J76").Value = 7173.8887
$ws.Range("K76").Value = 5378.5713
$ws.Range("L76").Value = 7173.8887
$ws.Range("M76").Value = -5063.5713
$ws.Range("N76").Value = -7803.8887
# Row 79
$ws.Range("H79").Value = 6388.4375
$ws.Range("I79").Value = 5378.5713
$ws.Range("J79").Value = 7173.8887
$ws.Range("K79").Value = 5378.5713
$ws.Range("L79").Value = 7173.8887
$ws.Range("M79").Value = -4286.5713
$ws.Range("N79").Value = -9357.8887
# Row 100
$ws.Range("H100").Value = 3010.182
$ws.Range("I100").Value = 2244.5715
$ws.Range("K100").Value = 2244.5715
$ws.Range("M100").Value = -1703.5715
# Row 113
$ws.Range("H113").Value = 3099.8572
$ws.Range("I113").Value = 2739.8
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 2739.8
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = 514.1999999999998
$ws.Range("N113").Value = -10508
# Row 124
$ws.Range("H124").Value = 35780
$ws.Range("J124").Value = 35780
$ws.Range("L124").Value = 35780
$ws.Range("N124").Value = -45600

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1826.9615
$ws.Range("I45").Value = 1763.25
$ws.Range("J45").Value = 1928.9
$ws.Range("K45").Value = 1763.25
$ws.Range("L45").Value = 1928.9
$ws.Range("M45").Value = -1386.25
$ws.Range("N45").Value = -2682.9
# Row 52
$ws.Range("H52").Value = 50780
$ws.Range("J52").Value = 50780
$ws.Range("L52").Value = 50780
$ws.Range("N52").Value = -51416
# Row 63
$ws.Range("H63").Value = 62502356
$ws.Range("I63").Value = 83335490
$ws.Range("J63").Value = 2965
$ws.Range("K63").Value = 83335490
$ws.Range("L63").Value = 2965
$ws.Range("M63").Value = -83334804
$ws.Range("N63").Value = -4337
# Row 66
$ws.Range("H66").Value = 62502356
$ws.Range("I66").Value = 83335490
$ws.Range("J66").Value = 2965
$ws.Range("K66").Value = 416677450
$ws.Range("L66").Value = 14825
$ws.Range("M66").Value = -416674018
$ws.Range("N66").Value = -21689
# Row 86
$ws.Range("H86").Value = 20000
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").Value = $null
# Row 89
$ws.Range("H89").Value = 20000
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").Value = $null
# Row 132
$ws.Range("H132").Value = 6758657.5
$ws.Range("I132").Value = 8930043
$ws.Range("J132").Value = 3234.6667
$ws.Range("K132").Value = 26790129
$ws.Range("L132").Value = 9704.000100000001
$ws.Range("M132").Value = -26787599
$ws.Range("N132").Value = -14764.0001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 4088.725
$ws.Range("I105").Value = 3279
$ws.Range("J105").Value = 4395.8623
$ws.Range("K105").Value = 3279
$ws.Range("L105").Value = 4395.8623
$ws.Range("M105").Value = -1532
$ws.Range("N105").Value = -7889.8623

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 14498694
$ws.Range("I31").Value = 6580.55
$ws.Range("J31").Value = 111112780
$ws.Range("K31").Value = 6580.55
$ws.Range("L31").Value = 111112780
$ws.Range("M31").Value = -6285.55
$ws.Range("N31").Value = -111113370
# Row 34
$ws.Range("H34").Value = 14498694
$ws.Range("I34").Value = 6580.55
$ws.Range("J34").Value = 111112780
$ws.Range("K34").Value = 6580.55
$ws.Range("L34").Value = 111112780
$ws.Range("M34").Value = -6378.55
$ws.Range("N34").Value = -111113184
# Row 99
$ws.Range("H99").Value = 1459.875
$ws.Range("I99").Value = 1090
$ws.Range("J99").Value = 1897
$ws.Range("K99").Value = 1090
$ws.Range("L99").Value = 1897
$ws.Range("M99").Value = 408
$ws.Range("N99").Value = -4893
# Row 105
$ws.Range("H105").Value = 1820.08
$ws.Range("I105").Value = 1154
$ws.Range("K105").Value = 1154
$ws.Range("M105").Value = 593
# Row 126
$ws.Range("H126").Value = 1459.875
$ws.Range("I126").Value = 1090
$ws.Range("J126").Value = 1897
$ws.Range("K126").Value = 3270
$ws.Range("L126").Value = 5691
$ws.Range("M126").Value = -800
$ws.Range("N126").Value = -10631
# Row 132
$ws.Range("H132").Value = 2293.375
$ws.Range("I132").Value = 1660.6774
$ws.Range("J132").Value = 4472.6665
$ws.Range("K132").Value = 4982.0322
$ws.Range("L132").Value = 13417.9995
$ws.Range("M132").Value = -2452.0322
$ws.Range("N132").Value = -18477.9995
# Row 140
$ws.Range("H140").Value = 38161.816
$ws.Range("J140").Value = 38161.816
$ws.Range("L140").Value = 38161.816
$ws.Range("N140").Value = -48521.816

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").Value = $null
# Row 102
$ws.Range("H102").Value = 2459.4375
$ws.Range("I102").Value = 2883.5454
$ws.Range("J102").Value = 1526.4
$ws.Range("K102").Value = 2883.5454
$ws.Range("L102").Value = 1526.4
$ws.Range("M102").Value = -1261.5454
$ws.Range("N102").Value = -4770.4
# Row 126
$ws.Range("H126").Value = 3865.5
$ws.Range("I126").Value = 2781.4546
$ws.Range("J126").Value = 4566.9414
$ws.Range("K126").Value = 8344.363799999999
$ws.Range("L126").Value = 13700.8242
$ws.Range("M126").Value = -5874.363799999999
$ws.Range("N126").Value = -18640.8242
# Row 132
$ws.Range("H132").Value = 5671.696
$ws.Range("I132").Value = 4766.3335
$ws.Range("J132").Value = 6253.7144
$ws.Range("K132").Value = 14299.0005
$ws.Range("L132").Value = 18761.1432
$ws.Range("M132").Value = -11769.0005
$ws.Range("N132").Value = -23821.1432
# Row 138
$ws.Range("H138").Value = 57582.332
$ws.Range("J138").Value = 57582.332
$ws.Range("L138").Value = 57582.332
$ws.Range("N138").Value = -67862.33199999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3131.195
$ws.Range("I40").Value = 4040.5264
$ws.Range("J40").Value = 2345.8635
$ws.Range("K40").Value = 4040.5264
$ws.Range("L40").Value = 2345.8635
$ws.Range("M40").Value = -3904.5264
$ws.Range("N40").Value = -2617.8635
# Row 122
$ws.Range("H122").Value = 7073.8184
$ws.Range("I122").Value = 7581.1665
$ws.Range("J122").Value = 6465
$ws.Range("K122").Value = 22743.4995
$ws.Range("L122").Value = 19395
$ws.Range("M122").Value = -20293.4995
$ws.Range("N122").Value = -24295
# Row 139
$ws.Range("H139").Value = 42025.4
$ws.Range("J139").Value = 42178.223
$ws.Range("L139").Value = 42178.223
$ws.Range("N139").Value = -52458.223

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2038.2963
$ws.Range("I122").Value = 1843.409
$ws.Range("J122").Value = 2895.8
$ws.Range("K122").Value = 5530.227000000001
$ws.Range("L122").Value = 8687.400000000001
$ws.Range("M122").Value = -3080.227000000001
$ws.Range("N122").Value = -13587.4

Write-Output "Applied all cell updates."